$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.051.81'
$ws.Range("E2").Value = '  -2.12%  '

$ws.Range("D3").Value = '2.496.83'
$ws.Range("E3").Value = '  -1.60%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '538.29'
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("D6").Value = '137.93'
$ws.Range("E6").Value = '  -5.16%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("E8").Value = '  -2.25%  '

$ws.Range("D9").Value = '2.494.24'
$ws.Range("E9").Value = '  -3.00%  '

$ws.Range("E10").Value = '  -1.45%  '

$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("D12").Value = '5.43'
$ws.Range("E12").Value = '  -2.05%  '

$ws.Range("D13").Value = '0.348'
$ws.Range("E13").Value = '  -4.20%  '

$ws.Range("D14").Value = '2.950.13'
$ws.Range("E14").Value = '  -1.09%  '

$ws.Range("D15").Value = '22.93'
$ws.Range("E15").Value = '  -5.41%  '

$ws.Range("D16").Value = '58.950.90'
$ws.Range("E16").Value = '  -2.17%  '

$ws.Range("E17").Value = '  -2.55%  '

$ws.Range("D18").Value = '2.495.44'
$ws.Range("E18").Value = '  -2.14%  '

$ws.Range("D19").Value = '10.94'
$ws.Range("E19").Value = '  -3.77%  '

$ws.Range("D20").Value = '4.25'
$ws.Range("E20").Value = '  -2.64%  '

$ws.Range("D21").Value = '323.96'
$ws.Range("E21").Value = '  -1.50%  '

$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").Value = '5.81'
$ws.Range("E23").Value = '  -1.95%  '

$ws.Range("D24").Value = '62.90'
$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("E25").Value = '  -6.02%  '

$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("E27").Value = '  +0.60%  '

$ws.Range("D28").Value = '7.64'
$ws.Range("E28").Value = '  -4.78%  '

$ws.Range("D29").Value = '0.0₃0771'
$ws.Range("E29").Value = '  -3.60%  '

$ws.Range("D30").Value = '6.59'
$ws.Range("E30").Value = '  -7.87%  '

$ws.Range("E31").Value = '  -1.97%  '

$ws.Range("D32").Value = '165.27'
$ws.Range("E32").Value = '  +1.61%  '

$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.21%  '

$ws.Range("E34").Value = '  -11.88%  '

$ws.Range("E35").Value = '  -8.32%  '

$ws.Range("E36").Value = '  -1.99%  '

$ws.Range("E37").Value = '  -9.42%  '

$ws.Range("E38").Value = '  -5.70%  '

$ws.Range("E39").Value = '  -3.46%  '

$ws.Range("D40").Value = '0.797'
$ws.Range("E40").Value = '  -5.30%  '

$ws.Range("D41").Value = '5.16'
$ws.Range("E41").Value = '  -9.72%  '

$ws.Range("D42").Value = '276.13'
$ws.Range("E42").Value = '  -9.50%  '

$ws.Range("E43").Value = '  +0.77%  '

$ws.Range("E44").Value = '  +0.51%  '

$ws.Range("D45").Value = '0.591'
$ws.Range("E45").Value = '  -2.84%  '

$ws.Range("D46").Value = '0.0935'
$ws.Range("E46").Value = '  -0.64%  '

$ws.Range("D47").Value = '124.46'
$ws.Range("E47").Value = '  -0.19%  '

$ws.Range("E48").Value = '  -3.22%  '

$ws.Range("E49").Value = '  -4.34%  '

$ws.Range("D50").Value = '17.49'
$ws.Range("E50").Value = '  -4.88%  '

$ws.Range("D51").Value = '1.764.84'
$ws.Range("E51").Value = '  -2.90%  '
